$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 13 (PPGEO/GEOGRAFIA's new "Sonda multiparametros"
# item), pushing the previous rows 13-19 down to 14-20.
$ws.Rows.Item(13).Insert()

$justificativaAgua = "O equipamento apoiará a realização de projetos que assumam a água como fator ambiental central de análise, em diálogo com condições pristinas e alterações ambientais em hidrossistemas; projetos que tenham na água o elemento-chave de compreensão da dinâmica das paisagens; projetos relacionados à poluição hídrica e saúde humana.`nProfessores`nMiguel Felippe - TERRA/PPGEO`nCézar Rocha - NAGEA/PPGEO/PROAC`nCamila Neves - GEOPED/PPGEO`n"

$ws.Cells.Item(13, 1).Value = "GEOGRAFIA"
$ws.Cells.Item(13, 2).Value = 2
$ws.Cells.Item(13, 3).Value = "Sonda multiparâmetros de qualidade de água"
$ws.Cells.Item(13, 4).Value = 1
$ws.Cells.Item(13, 5).Value = 8000
$ws.Cells.Item(13, 6).Value = 8000
$ws.Cells.Item(13, 7).Value = "P1 - Alta"
$ws.Cells.Item(13, 8).Value = $justificativaAgua
$ws.Cells.Item(13, 9).Value = "GEOGRAFIA, AMBIENTE CONSTRUÍDO"

# Writing a multi-line justification triggers Excel's row auto-height; reset
# it back to the sheet's default so row 13 matches the other rows.
$ws.Rows.Item(13).EntireRow.AutoFit()
